# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.621.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = "'2.552.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'309.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.09%  '
$ws.Range("D6").Value = "'97.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("D7").Value = "'0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.21%  '
$ws.Range("D9").Value = "'0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("D10").Value = "'35.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.15%  '
$ws.Range("D11").Value = "'0.0804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = "'7.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.30%  '
$ws.Range("D13").Value = "'2.949.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").Value = "'15.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.68%  '
$ws.Range("D16").Value = "'2.611.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.53%  '
$ws.Range("D17").Value = "'0.832"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").Value = "'42.673.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").Value = "'6.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.65%  '
$ws.Range("D20").Value = "'0.0₃0952"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Value = "'12.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.88%  '
$ws.Range("D22").Value = "'69.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.67%  '
$ws.Range("D23").Value = "'246.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("D24").Value = "'2.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = "'2.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").Value = "'26.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = "'2.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("D29").Value = "'39.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.33%  '
$ws.Range("D30").Value = "'10.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.54%  '
$ws.Range("D31").Value = "'157.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.06%  '
$ws.Range("D32").Value = "'5.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.88%  '
$ws.Range("D33").Value = "'0.0790"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").Value = "'3.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.06%  '
$ws.Range("D35").Value = "'2.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.98%  '
$ws.Range("D36").Value = "'2.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.13%  '
$ws.Range("D37").Value = "'18.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("E38").Value = '  +7.03%  '
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").Value = "'22.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("D42").Value = "'4.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.59%  '
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = "'0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.75%  '
$ws.Range("D45").Value = "'1.983.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.05%  '
$ws.Range("D46").Value = "'3.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.95%  '
$ws.Range("D47").Value = "'8.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("D48").Value = "'2.802.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").Value = "'80.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.70%  '
$ws.Range("D50").Value = "'0.191"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").Value = "'72.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.18%  '
